# Edits to SamMcGrailResume-2022.docx
# Each block below merges runs that were previously split by
# <w:proofErr/> spell/grammar markers (or tab-separated fragments) back
# into single runs, and makes the couple of small text changes
# (Mar -> Apr, and the "Oct" spacing tweak).

$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceAll = 2

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, $wdFindContinue, $false, $replace, $wdReplaceAll)
}

# 1) " github.com/" + "sammcgrail" + " " -> " github.com/sammcgrail "
Replace-Text " github.com/sammcgrail " " github.com/sammcgrail "

# 2) "linkedin.com/in/" + "sammcgrail" + " " -> "linkedin.com/in/sammcgrail "
Replace-Text "linkedin.com/in/sammcgrail " "linkedin.com/in/sammcgrail "

# 3) "Glaukos" + " " -> "Glaukos "
Replace-Text "Glaukos " "Glaukos "

# 4) "doblePRIME" + "  " -> "doblePRIME  "
Replace-Text "doblePRIME  " "doblePRIME  "

# 5) "Extensive use of various " + "javascript" -> merged run
Replace-Text "Extensive use of various javascript" "Extensive use of various javascript"

# 6) "Launch " + "Academy  " -> "Launch Academy  "
Replace-Text "Launch Academy  " "Launch Academy  "

# 7) "Junior Full Stack Web " + "Developer  " -> merged run
Replace-Text "Junior Full Stack Web Developer  " "Junior Full Stack Web Developer  "

# 8) "Professional " + "Musician" -> "Professional Musician"
Replace-Text "Professional Musician" "Professional Musician"

# 9) ": React, Webpack, Node, *." + "js" + ", shell-fu..." -> merged run
Replace-Text ": React, Webpack, Node, *.js, shell-fu, git, bash, custom VM curation" `
             ": React, Webpack, Node, *.js, shell-fu, git, bash, custom VM curation"

# 10) Date fix: Mar -> Apr (the "Mar 2022 (Current)" occurrence)
Replace-Text "Mar 2022 (Current)" "Apr 2022 (Current)"
